# Updated cryptos list values (Price / Volume(1h) columns, and a couple of row
# swaps where the source ranking reordered two coins) to match the refreshed
# data snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    # Force the cell to stay a text value even when it looks like a number
    # (e.g. "229.74"), mirroring how Excel treats a leading apostrophe on
    # manual entry, then drop the resulting quote-prefix style so the cell
    # keeps its original (default) formatting.
    $Range.Value = "'" + $Text
    $Range.Style = "Normal"
}

$ws.Range('D2').Value = '37.409.59'
$ws.Range('E2').Value = '  -1.05%  '
$ws.Range('D3').Value = '2.051.04'
$ws.Range('E3').Value = '  -1.53%  '
Set-TextValue $ws.Range('D4') '0.999'
$ws.Range('E4').Value = '  -0.15%  '
Set-TextValue $ws.Range('D5') '229.74'
$ws.Range('E5').Value = '  -1.60%  '
Set-TextValue $ws.Range('D6') '0.615'
$ws.Range('E6').Value = '  -1.59%  '
Set-TextValue $ws.Range('D8') '57.33'
$ws.Range('E8').Value = '  -2.08%  '
$ws.Range('E9').Value = '  -1.66%  '
Set-TextValue $ws.Range('D10') '0.0798'
$ws.Range('E10').Value = '  +1.41%  '
$ws.Range('E11').Value = '  -2.08%  '
Set-TextValue $ws.Range('D12') '14.81'
$ws.Range('E12').Value = '  -0.73%  '
$ws.Range('D13').Value = '2.354.66'
$ws.Range('E13').Value = '  -1.47%  '
Set-TextValue $ws.Range('D14') '20.81'
$ws.Range('E14').Value = '  -1.73%  '
$ws.Range('E15').Value = '  -3.07%  '
Set-TextValue $ws.Range('D16') '5.32'
$ws.Range('E16').Value = '  -0.73%  '
$ws.Range('D17').Value = '2.056.02'
$ws.Range('E17').Value = '  -1.75%  '
$ws.Range('D18').Value = '37.312.76'
$ws.Range('E18').Value = '  -1.24%  '
Set-TextValue $ws.Range('D19') '6.08'
$ws.Range('E19').Value = '  -1.28%  '
Set-TextValue $ws.Range('D20') '69.70'
$ws.Range('E20').Value = '  -2.58%  '
$ws.Range('D21').Value = '0.0₃0833'
$ws.Range('E21').Value = '  -0.97%  '
Set-TextValue $ws.Range('D22') '226.66'
$ws.Range('E22').Value = '  -1.23%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('E24').Value = '  -0.35%  '
Set-TextValue $ws.Range('D25') '2.30'
$ws.Range('E25').Value = '  -3.93%  '
$ws.Range('E26').Value = '  -2.38%  '
Set-TextValue $ws.Range('D27') '168.92'
$ws.Range('E27').Value = '  -1.73%  '
$ws.Range('E28').Value = '  -6.21%  '
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D29') '1.38'
$ws.Range('E29').Value = '  -1.61%  '
$ws.Range('B30').Value = 'EthereumClassic'
$ws.Range('C30').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D30') '19.02'
$ws.Range('E30').Value = '  -2.60%  '
$ws.Range('E31').Value = '  -2.22%  '
Set-TextValue $ws.Range('D32') '4.57'
$ws.Range('E32').Value = '  -3.54%  '
$ws.Range('E33').Value = '  -2.49%  '
$ws.Range('E34').Value = '  -1.45%  '
Set-TextValue $ws.Range('D35') '2.45'
$ws.Range('E35').Value = '  -0.23%  '
Set-TextValue $ws.Range('D36') '1.84'
$ws.Range('E36').Value = '  +1.46%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D37') '3.26'
$ws.Range('E37').Value = '  -4.11%  '
$ws.Range('B38').Value = 'BinanceUSD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue $ws.Range('D38') '1.00'
$ws.Range('E38').Value = '  -0.14%  '
Set-TextValue $ws.Range('D39') '5.36'
$ws.Range('E39').Value = '  -1.34%  '
$ws.Range('E40').Value = '  -5.14%  '
Set-TextValue $ws.Range('D41') '17.27'
$ws.Range('E41').Value = '  +2.42%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.485.19'
$ws.Range('E42').Value = '  +2.73%  '
$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D43') '2.89'
$ws.Range('E43').Value = '  -1.14%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D44') '0.0944'
$ws.Range('E44').Value = '  -3.24%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D45') '96.61'
$ws.Range('E45').Value = '  -4.58%  '
$ws.Range('E46').Value = '  +0.91%  '
$ws.Range('E47').Value = '  -4.19%  '
Set-TextValue $ws.Range('D48') '3.95'
$ws.Range('E48').Value = '  -3.76%  '
$ws.Range('E49').Value = '  -2.32%  '
Set-TextValue $ws.Range('D50') '2.93'
$ws.Range('E50').Value = '  -2.03%  '
$ws.Range('D51').Value = '2.243.37'
$ws.Range('E51').Value = '  -1.42%  '
